{"js": "// Update the date line and the 25 division-fact table cells to the new\n// values generated for this day's worksheet.\n\n// --- Title paragraph: date string -----------------------------------\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\n\n// --- Table cells ------------------------------------------------------\n// Rows 0,4,8,12,16 (0-based) hold the visible answers; the rows between\n// them are spacer rows. Each holds 5 columns of \"a\u00f7b=c, r\" text.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\n\nawait context.sync();\n\nconst titlePar = paragraphs.items[0];\nconst table = tables.items[0];\n\nconst newValues = [\n  [\"33\u00f75=6, 3\", \"71\u00f78=8, 7\", \"23\u00f74=5, 3\", \"34\u00f76=5, 4\", \"61\u00f73=20, 1\"],\n  [\"50\u00f78=6, 2\", \"71\u00f78=8, 7\", \"29\u00f75=5, 4\", \"13\u00f78=1, 5\", \"10\u00f76=1, 4\"],\n  [\"87\u00f76=14, 3\", \"21\u00f76=3, 3\", \"47\u00f78=5, 7\", \"94\u00f76=15, 4\", \"63\u00f72=31, 1\"],\n  [\"76\u00f77=10, 6\", \"97\u00f73=32, 1\", \"57\u00f73=19, 0\", \"88\u00f77=12, 4\", \"66\u00f76=11, 0\"],\n  [\"37\u00f78=4, 5\", \"43\u00f76=7, 1\", \"54\u00f78=6, 6\", \"96\u00f74=24, 0\", \"14\u00f76=2, 2\"],\n];\nconst dataRows = [0, 4, 8, 12, 16];\n\n// Gather the (single) paragraph inside every target cell so we can\n// replace its range's text without disturbing run formatting (font,\n// size, paragraph alignment, etc.)\nconst cellParagraphCollections = [];\nfor (let i = 0; i < dataRows.length; i++) {\n  const rowIndex = dataRows[i];\n  for (let col = 0; col < 5; col++) {\n    const cell = table.getCell(rowIndex, col);\n    const cellParagraphs = cell.body.paragraphs;\n    cellParagraphs.load(\"items\");\n    cellParagraphCollections.push(cellParagraphs);\n  }\n}\n\nawait context.sync();\n\n// Title paragraph: replace its text in place (keeps rFonts/sz/jc).\ntitlePar.getRange().insertText(\"2024-09-04 Wednesday\", Word.InsertLocation.replace);\n\n// Table cells: same in-place range replacement, one per cell.\nlet idx = 0;\nfor (let i = 0; i < dataRows.length; i++) {\n  for (let col = 0; col < 5; col++) {\n    const cellParagraphs = cellParagraphCollections[idx++];\n    const cellRange = cellParagraphs.items[0].getRange();\n    cellRange.insertText(newValues[i][col], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division-fact table cells to the new\n# values generated for this day's worksheet.\n\n$d = $word.ActiveDocument\n\n# --- Title paragraph: date string -----------------------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"2024-09-03 Tuesday\"\n$find.Replacement.Text = \"2024-09-04 Wednesday\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# --- Table cells ------------------------------------------------------\n# Rows 1,5,9,13,17 (1-based) hold the visible answers; the rows between\n# them are spacer rows. Each holds 5 columns of \"a\u00f7b=c, r\" text.\n$t = $d.Tables.Item(1)\n\n$newValues = @{\n    \"1,1\" = \"33\u00f75=6, 3\";\n    \"1,2\" = \"71\u00f78=8, 7\";\n    \"1,3\" = \"23\u00f74=5, 3\";\n    \"1,4\" = \"34\u00f76=5, 4\";\n    \"1,5\" = \"61\u00f73=20, 1\";\n    \"5,1\" = \"50\u00f78=6, 2\";\n    \"5,2\" = \"71\u00f78=8, 7\";\n    \"5,3\" = \"29\u00f75=5, 4\";\n    \"5,4\" = \"13\u00f78=1, 5\";\n    \"5,5\" = \"10\u00f76=1, 4\";\n    \"9,1\" = \"87\u00f76=14, 3\";\n    \"9,2\" = \"21\u00f76=3, 3\";\n    \"9,3\" = \"47\u00f78=5, 7\";\n    \"9,4\" = \"94\u00f76=15, 4\";\n    \"9,5\" = \"63\u00f72=31, 1\";\n    \"13,1\" = \"76\u00f77=10, 6\";\n    \"13,2\" = \"97\u00f73=32, 1\";\n    \"13,3\" = \"57\u00f73=19, 0\";\n    \"13,4\" = \"88\u00f77=12, 4\";\n    \"13,5\" = \"66\u00f76=11, 0\";\n    \"17,1\" = \"37\u00f78=4, 5\";\n    \"17,2\" = \"43\u00f76=7, 1\";\n    \"17,3\" = \"54\u00f78=6, 6\";\n    \"17,4\" = \"96\u00f74=24, 0\";\n    \"17,5\" = \"14\u00f76=2, 2\";\n}\n\n$rows = @(1, 5, 9, 13, 17)\nforeach ($r in $rows) {\n    for ($c = 1; $c -le 5; $c++) {\n        $key = \"$r,$c\"\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newValues[$key]\n    }\n}\n"}
